# feat: add kuota field to matkul and update related features
#
# Adds two new columns (kode_matkul, semester) and a kuota column's worth of
# data to the existing "matkul" header/data row, turning:
#   nama | kode
#   coba | coba
# into:
#   nama | kuota | kode_matkul | semester
#   coba | 1     | coba        | 1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename B1 ("kode" -> "kuota") and add the two new headers.
$ws.Range("B1").Value = "kuota"
$ws.Range("C1").Value = "kode_matkul"
$ws.Range("D1").Value = "semester"

# Data row: B2 becomes a numeric kuota value, C2 mirrors the "coba" sample
# text, D2 is a numeric semester value. A2 keeps its "coba" text.
$ws.Range("A2").Value = "coba"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "coba"
$ws.Range("D2").Value = 1

# Match the new active selection recorded in the saved file.
[void]$ws.Range("E3").Select()
